$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H3").Value = 3.8
$ws.Range("J3").Value = 1.03
$ws.Range("K3").Value = 15
$ws.Range("L3").Value = 1.18
$ws.Range("M3").Value = 4.5
$ws.Range("N3").Value = 1.62
$ws.Range("O3").Value = 2.25
$ws.Range("P3").Value = 1.33
$ws.Range("Q3").Value = 3.25
$ws.Range("R3").Value = 1.67
$ws.Range("S3").Value = 2.1
$ws.Range("T3").Value = 8.5
$ws.Range("U3").Value = 8.5
$ws.Range("Z3").Value = 13
$ws.Range("AA3").Value = 7.5
$ws.Range("AC3").Value = 41
$ws.Range("AD3").Value = 151
$ws.Range("AG3").Value = 15
$ws.Range("AI3").Value = 34
$ws.Range("K5").Value = 17
$ws.Range("L5").Value = 1.14
$ws.Range("M5").Value = 5.5
$ws.Range("N5").Value = 1.5
$ws.Range("O5").Value = 2.5
$ws.Range("P5").Value = 1.29
$ws.Range("Q5").Value = 3.5
$ws.Range("R5").Value = 3
$ws.Range("S5").Value = 1.36
$ws.Range("T5").Value = 7
$ws.Range("U5").Value = 5.5
$ws.Range("V5").Value = 12
$ws.Range("X5").Value = 13
$ws.Range("Y5").Value = 41
$ws.Range("Z5").Value = 13
$ws.Range("AC5").Value = 151
$ws.Range("AI5").Value = 201
$ws.Range("AJ5").Value = 151
$ws.Range("G7").Value = 2.85
$ws.Range("H7").Value = 2.72
$ws.Range("I7").Value = 2.7
$ws.Range("W7").Value = 37
$ws.Range("X7").Value = 25
$ws.Range("Y7").Value = 32
$ws.Range("Z7").Value = 7.3
$ws.Range("AE7").Value = 7.6
$ws.Range("AG7").Value = 9.75
$ws.Range("AH7").Value = 35
$ws.Range("AI7").Value = 25
$ws.Range("G8").Value = 5.9
$ws.Range("H8").Value = 4.35
$ws.Range("I8").Value = 1.4
$ws.Range("O8").Value = 2.1
$ws.Range("R8").Value = 1.81
$ws.Range("S8").Value = 1.9
$ws.Range("T8").Value = 15
$ws.Range("U8").Value = 30
$ws.Range("V8").Value = 16
$ws.Range("W8").Value = 90
$ws.Range("X8").Value = 45
$ws.Range("Y8").Value = 45
$ws.Range("AA8").Value = 7.7
$ws.Range("AB8").Value = 14.5
$ws.Range("AC8").Value = 55
$ws.Range("AE8").Value = 6.7
$ws.Range("AF8").Value = 6.2
$ws.Range("AH8").Value = 8
$ws.Range("AI8").Value = 9.25
$ws.Range("AJ8").Value = 19
$ws.Range("N10").Value = 1.48
$ws.Range("O10").Value = 2.6
$ws.Range("Z10").Value = 19
$ws.Range("G11").Value = 2.35
$ws.Range("I11").Value = 2.88
$ws.Range("U11").Value = 13
$ws.Range("V11").Value = 9.5
$ws.Range("AI11").Value = 21
